# Update "want to go" (想去人数) counts scraped at a later point in time.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 1451
$ws1.Range("F5").Value = 188
$ws1.Range("F6").Value = 38
$ws1.Range("F8").Value = 9724
$ws1.Range("F9").Value = 164
$ws1.Range("F10").Value = 110
$ws1.Range("F13").Value = 368
$ws1.Range("F14").Value = 6709
$ws1.Range("F15").Value = 1081
$ws1.Range("F16").Value = 124
$ws1.Range("F17").Value = 50
$ws1.Range("F18").Value = 184

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 262

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 1451
$ws4.Range("F5").Value = 188
$ws4.Range("F6").Value = 38
$ws4.Range("F8").Value = 262
$ws4.Range("F10").Value = 9724
$ws4.Range("F11").Value = 164
$ws4.Range("F12").Value = 110
$ws4.Range("F15").Value = 368
$ws4.Range("F16").Value = 6709
$ws4.Range("F17").Value = 1081
$ws4.Range("F18").Value = 124
$ws4.Range("F19").Value = 50
$ws4.Range("F20").Value = 184
